# "I think I got an idea"
#
# 1. "Emily is Away + ILOVEYOU virus" -> "Emily is Away +" / " bug in code"
# 2. Brand new paragraph right after it with the full game-idea brainstorm.

$d = $word.ActiveDocument

# --- 1. Swap "ILOVEYOU virus" for "bug in code" on the second paragraph ---
$bugPara = $d.Paragraphs(2).Range
$hit = $d.Range($bugPara.Start, $bugPara.End)
$found = $hit.Find.Execute(" ILOVEYOU virus")
if ($found -and $hit.Find.Found) {
    $hit.Text = ""
}

$bugPara = $d.Paragraphs(2).Range
$tail = $d.Range($bugPara.End - 1, $bugPara.End - 1)
$tail.InsertAfter(" bug in code")

# --- 2. New paragraph right after it with the brainstorm dump ---
$bugPara = $d.Paragraphs(2).Range
$bugPara.InsertParagraphAfter()

$apos = [char]0x2019   # RIGHT SINGLE QUOTATION MARK
$endash = [char]0x2013 # EN DASH

$ideaText = "Maybe main character is talking to someone (love interest) after a few days, " + `
    "goes missing, cortana/siri kinda assistant offers to help (so talking like future stuff " + `
    "here). After some time the ai starts to get some stuff wrong. Maybe AI assistant has like " + `
    "personality sliders in settings or something? And they change over the course of the game " + `
    "(player can change them in settings, but they don" + $apos + "t save). IVE GOT IT " + $endash + `
    " INTENTIONAL BUGS " + $endash + " MALICIOUS PERSON!!!! "

$ideaPara = $d.Paragraphs(3).Range
$ideaStart = $d.Range($ideaPara.Start, $ideaPara.Start)
$ideaStart.InsertAfter($ideaText)

Write-Output $d.Paragraphs(2).Range.Text
Write-Output $d.Paragraphs(3).Range.Text
Write-Output ("ParagraphCount=" + $d.Paragraphs.Count)
